$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1 (22:22 -> 22:52)
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 22:52"

# Estados Unidos (row 4) - updated covid figures
$ws.Range("B4").Value = 495750
$ws.Range("C4").Value = 27184
$ws.Range("E4").Value = 450537
$ws.Range("G4").Value = 1739
$ws.Range("H4").Value = 18430

# Alemania (row 8) - updated covid figures
$ws.Range("B8").Value = 121045
$ws.Range("C8").Value = 2810
$ws.Range("E8").Value = 64404
$ws.Range("G8").Value = 121
$ws.Range("H8").Value = 2728

# Rows 129/130: El Salvador and Ruanda swap places (Ruanda now listed before
# El Salvador), each carrying its own updated data.
$ws.Range("A129").Value = "Ruanda"
$ws.Range("B129").Value = 118
$ws.Range("C129").Value = 5
$ws.Range("D129").Value = 7
$ws.Range("E129").Value = 111
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 0

$ws.Range("A130").Value = "El Salvador"
$ws.Range("B130").Value = 117
$ws.Range("C130").Value = 14
$ws.Range("D130").Value = 15
$ws.Range("E130").Value = 96
$ws.Range("F130").Value = 4
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 6
